$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, pushing the current rows 2-5 down to 3-6.
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the new observation record.
$ws.Range("A2").Value = 95503476
$ws.Range("B2").Value = 96334
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."

# "Antal" is stored as text in this sheet (e.g. "30", "15"), so force text
# formatting before assigning, otherwise Excel auto-converts the numeric
# looking string to a number.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "5"

$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("K2").Value = "överblommad"
$ws.Range("P2").Value = "A37438, Troserum, V Ed, Sm"
$ws.Range("Q2").Value = 591590.0879378035
$ws.Range("R2").Value = 6428792.104066458
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Kalmar"
$ws.Range("U2").Value = "Västervik"
$ws.Range("V2").Value = "Småland"
$ws.Range("W2").Value = "Västra Ed"

# Start/end dates are stored as plain text ("YYYY-MM-DD"), not real date
# values, so force text formatting before assigning to stop Excel from
# auto-converting them into date serials.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2021-08-15"

$ws.Range("Z2").Value = "10:30"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2021-08-15"

$ws.Range("AB2").Value = "14:00"
$ws.Range("AC2").Value = "1 blomma"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = "Magnus Kasselstrand"
$ws.Range("AX2").Value = "Magnus Kasselstrand, Ingvor Kasselstrand"
